$d = $word.ActiveDocument

function Replace-InRange($range, [string]$old, [string]$new) {
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. Update activation date
Replace-InRange $d.Content "Ativação: 01/01/2021" "Ativação: 01/01/2024"

# 2. "Programa resumido" - Portuguese summary paragraph
$pResumidoPt = $d.Paragraphs(11).Range
Replace-InRange $pResumidoPt `
    "1. Conceitos e Definições no Estudo da Inovação. 2. Tipos de inovação. 3. Ambientes de inovação. 4. Sistemas de Inovação. 5. Capacidade tecnológica e aprendizagem. 6. Processos de inovação. 7. Estratégias de inovação. 8. Políticas públicas para promoção da inovação. 9. Propriedade Intelectual." `
    "1. Conceitos e Definições no Estudo da Inovação. 2. Tipos de inovação. 3. Fontes de inovação; 4. Sistemas de Inovação. 5. Capacidade tecnológica e aprendizagem. 6. Processos de inovação. 7. Estratégias de inovação. 8. Cultura organizacional e inovação. 9. Comercialização e difusão da inovação."

# 3. "Programa resumido" - English summary paragraph (italic)
$pResumidoEn = $d.Paragraphs(12).Range
Replace-InRange $pResumidoEn `
    "1. Innovation: definition and key concepts. 2. Types of innovation. 3. Innovation environments. 4. Innovation System. 5. Technological Capabilities and Learning. 6. Innovation process. 7. Innovation strategy. 8. Public policies to promote innovation. 9. Intellectual Property." `
    "1. Innovation: definition and key concepts. 2. Types of innovation. Sources of innovation. 4. Innovation System. 5. Technological Capabilities and Learning. 6. Innovation process. 7. Innovation strategy.  8. Organizacional culture and innovation. 9.Commercialization and diffusion of innovations."

# 4. "Programa" - Portuguese full program paragraph
$pProgramaPt = $d.Paragraphs(14).Range
Replace-InRange $pProgramaPt `
    ". Conceitos e Definições no Estudo da Inovação. 2. Tipos de inovação. 3. Ambientes de inovação. 4. Sistemas de Inovação. 5. Capacidade tecnológica e aprendizagem. 6. Processos de inovação. 7. Estratégias de inovação. 8. Políticas públicas para promoção da inovação. 9. Propriedade Intelectual." `
    "1.Conceitos e Definições no Estudo da Inovação. 2. Tipos de inovação. 3. Ambientes de inovação. 4. Sistemas de Inovação. 5. Capacidade tecnológica e aprendizagem. 6. Processos de inovação. 7. Estratégias de inovação; 8. Cultura organizacional e inovação. 9. Comercialização e difusão da inovação. 10. Estudos de casos."

# 5. "Programa" - English full program paragraph (italic)
$pProgramaEn = $d.Paragraphs(15).Range
Replace-InRange $pProgramaEn `
    "1. Innovation: definition and key concepts. 2. Types of innovation. 3. Innovation environments. 4. Innovation System. 5. Technological Capabilities and Learning. 6. Innovation process. 7. Innovation strategy. 8. Public policies to promote innovation. 9. Intellectual Property." `
    "1. Innovation: definition and key concepts. 2. Types of innovation. 3. Innovation environments. 4. Innovation System. 5. Technological Capabilities and Learning. 6. Innovation process. 7. Innovation strategy.  8. Organizacional culture and innovation. 9.Commercialization and diffusion of innovations. 10. Cases studies."

# 6. Avaliação - Método
$pAvaliacao = $d.Paragraphs(17).Range
Replace-InRange $pAvaliacao "Provas e Trabalhos" "Provas, atividades em grupo e atividades individuais."

# 7. Avaliação - Critério
Replace-InRange $pAvaliacao `
    "M = (0,8P + 0,2T)P = média aritmética de duas provas escritasT = Média das notas de trabalhos e exercíciosM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas." `
    "Média das atividades avaliativas"

# 8. Bibliografia - edition update
Replace-InRange $d.Content "2nd ed. Palgrave – MacMillan, Houndsmill, 2010." "3nd ed. Palgrave – MacMillan, Houndsmill, 2017."
